$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147, shifting the existing rows 147-157 down to 148-158.
$ws.Rows("147").Insert()

# Populate the newly inserted row 147 with the new weekly price record.
$ws.Cells.Item(147, 1).Value  = 5
$ws.Cells.Item(147, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(147, 3).Value  = "Maule"
$ws.Cells.Item(147, 4).Value  = 45008
$ws.Cells.Item(147, 5).Value  = 7
$ws.Cells.Item(147, 6).Value  = 100112001
$ws.Cells.Item(147, 7).Value  = "Berenjena"
$ws.Cells.Item(147, 8).Value  = "Sin especificar"
$ws.Cells.Item(147, 9).Value  = "Primera"
$ws.Cells.Item(147, 10).Value = 200
$ws.Cells.Item(147, 11).Value = 8000
$ws.Cells.Item(147, 12).Value = 8000
$ws.Cells.Item(147, 13).Value = 8000
$ws.Cells.Item(147, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(147, 15).Value = "Región del Maule"
$ws.Cells.Item(147, 16).Value = 160
$ws.Cells.Item(147, 17).Value = 50
$ws.Cells.Item(147, 18).Value = "Hortaliza"
